$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below hold plain numeric-looking text in the original workbook
# (e.g. "243.61", "0.0000249") stored as inline strings, not numbers.
# Force the Text number format on them before writing the new values so
# Excel keeps storing them as text instead of auto-converting to numbers
# (which would drop trailing zeros / change precision).
$textNumberFormatCells = @(
    "D5",
    "D6",
    "D7",
    "D8",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D19",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D27",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D46",
    "D47",
    "D48",
    "D49",
    "D50",
    "D51"
)
foreach ($cellRef in $textNumberFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the refreshed coin values (price / volume / and the two swapped rows).
$ws.Range("D2").Value = "91.455.23"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "3.104.69"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "243.61"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Value = "615.58"
$ws.Range("E6").Value = "  -1.60%  "
$ws.Range("D7").Value = "1.09"
$ws.Range("E7").Value = "  -4.70%  "
$ws.Range("D8").Value = "0.383"
$ws.Range("E8").Value = "  +2.59%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "3.104.30"
$ws.Range("E10").Value = "  +13.76%  "
$ws.Range("D11").Value = "0.742"
$ws.Range("E11").Value = "  -4.13%  "
$ws.Range("D12").Value = "0.205"
$ws.Range("E12").Value = "  +0.17%  "
$ws.Range("D13").Value = "0.0000249"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").Value = "5.61"
$ws.Range("E14").Value = "  +2.27%  "
$ws.Range("D15").Value = "34.53"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").Value = "91.504.19"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "3.683.28"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "3.114.89"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").Value = "3.72"
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "14.75"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").Value = "5.76"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "447.06"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").Value = "9.27"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").Value = "0.0000201"
$ws.Range("E24").Value = "  -7.81%  "
$ws.Range("D25").Value = "5.62"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "88.96"
$ws.Range("E26").Value = "  -3.33%  "
$ws.Range("D27").Value = "11.64"
$ws.Range("E27").Value = "  -2.64%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  +19.76%  "
$ws.Range("D31").Value = "0.227"
$ws.Range("E31").Value = "  -6.02%  "
$ws.Range("D32").Value = "0.167"
$ws.Range("E32").Value = "  -10.27%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").Value = "0.175"
$ws.Range("E33").Value = "  +3.38%  "
$ws.Range("D34").Value = "9.25"
$ws.Range("E34").Value = "  -1.26%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").Value = "0.987"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("D36").Value = "7.66"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "26.16"
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("D38").Value = "1.94"
$ws.Range("E38").Value = "  +1.15%  "
$ws.Range("D39").Value = "3.94"
$ws.Range("E39").Value = "  -3.72%  "
$ws.Range("D40").Value = "487.04"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("D41").Value = "1.30"
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Value = "0.433"
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("D43").Value = "3.42"
$ws.Range("E43").Value = "  -5.61%  "
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "157.55"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "0.695"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "1.89"
$ws.Range("E48").Value = "  -1.68%  "
$ws.Range("D49").Value = "1.33"
$ws.Range("E49").Value = "  -1.04%  "
$ws.Range("D50").Value = "44.08"
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("D51").Value = "4.37"
$ws.Range("E51").Value = "  -3.79%  "
